$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (match styling of the existing header row, e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows: (row, I, J)
$data = @(
    @(2, 1, 6),
    @(3, 1, 5),
    @(4, 1, 5),
    @(5, 1, 5),
    @(6, 1, 7),
    @(7, 1, 4),
    @(8, 1, 4),
    @(9, 1, 6),
    @(10, 1, 6),
    @(11, 1, 6),
    @(12, 1, 4),
    @(13, 1, 6),
    @(14, 1, 4),
    @(15, 6, 8),
    @(16, 1, 2),
    @(17, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
